{"js": "// Fix a handful of Serbian-Latin typos (\"slovne greske\") in the project\n// report text, and tag the (empty) first title paragraph's mark with the\n// sr-Latn-RS language, matching the document's other paragraphs.\n\n// 1) The very first paragraph in the body is an empty \"Title\" styled\n//    paragraph. Stamp its paragraph mark with the sr-Latn-RS language\n//    (this shows up as <w:rPr><w:lang w:val=\"sr-Latn-RS\"/></w:rPr> inside\n//    that paragraph's <w:pPr>).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\nconst firstParagraphMark = firstParagraph.getRange(\"Whole\");\nfirstParagraphMark.languageId = \"sr-Latn-RS\";\nawait context.sync();\n\n// 2) \"ucestano\u0161cu\" -> \"u\u010destano\u0161cu\" (missing caron on the \u010d).\nlet found = context.document.body.search(\"ucestano\u0161cu\", { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\nif (found.items.length > 0) {\n  found.items[0].insertText(\"u\u010destano\u0161cu\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3) \"10 sample/s\" -> \"10 odbiraka po sekundi\" (use proper Serbian wording\n//    instead of the stray English \"sample/s\").\nfound = context.document.body.search(\"10 sample/s\", { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\nif (found.items.length > 0) {\n  found.items[0].insertText(\"10 odbiraka po sekundi\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 4) \"du\u017eina 100 odbiraka\" -> \"du\u017eine 100 odbiraka\" (grammar fix: genitive).\nfound = context.document.body.search(\"du\u017eina 100\", { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\nif (found.items.length > 0) {\n  found.items[0].insertText(\"du\u017eine 100\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Fix a handful of Serbian-Latin typos (\"slovne greske\") in the project\n# report text, and tag the (empty) first title paragraph's mark with the\n# sr-Latn-RS language, matching the document's other paragraphs.\n\n$d = $word.ActiveDocument\n\n# 1) The very first paragraph in the document is an empty \"Title\" styled\n#    paragraph. Stamp its paragraph mark with the sr-Latn-RS language\n#    (this shows up as <w:rPr><w:lang w:val=\"sr-Latn-RS\"/></w:rPr> inside\n#    that paragraph's <w:pPr>).\n$firstParagraph = $d.Paragraphs.Item(1)\n$firstParagraph.Range.LanguageID = \"sr-Latn-RS\"\n\n# 2) \"ucestano\u0161cu\" -> \"u\u010destano\u0161cu\" (missing caron on the \u010d).\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"ucestano\u0161cu\", $false, $false, $false, $false, $false, $true, 1, $false, \"u\u010destano\u0161cu\", 2) | Out-Null\n\n# 3) \"10 sample/s\" -> \"10 odbiraka po sekundi\" (use proper Serbian wording\n#    instead of the stray English \"sample/s\").\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"10 sample/s\", $false, $false, $false, $false, $false, $true, 1, $false, \"10 odbiraka po sekundi\", 2) | Out-Null\n\n# 4) \"du\u017eina 100 odbiraka\" -> \"du\u017eine 100 odbiraka\" (grammar fix: genitive).\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"du\u017eina 100\", $false, $false, $false, $false, $false, $true, 1, $false, \"du\u017eine 100\", 2) | Out-Null\n"}
